$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.820.24'
$ws.Range("E2").Value = '  -2.26%  '
$ws.Range("D3").Value = '1.775.75'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  -0.23%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.83'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.79%  '
$ws.Range("E6").Value = '  -1.75%  '
$ws.Range("E7").Value = '  -0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '31.11'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -4.46%  '
$ws.Range("E9").Value = '  +0.56%  '
$ws.Range("E10").Value = '  +5.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0921'
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = '2.029.64'
$ws.Range("E12").Value = '  -0.85%  '
$ws.Range("D13").Value = '1.772.61'
$ws.Range("E13").Value = '  -1.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.51'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.625'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.24%  '
$ws.Range("D16").Value = '33.817.24'
$ws.Range("E16").Value = '  -2.33%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '4.21'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '67.82'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.51%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.08'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.96%  '
$ws.Range("D20").Value = '0.0₃0774'
$ws.Range("E20").Value = '  +1.03%  '
$ws.Range("E21").Value = '  -0.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.58'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.72%  '
$ws.Range("E23").Value = '  -3.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.32'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '16.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.22%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.77%  '
$ws.Range("E28").Value = '  -2.10%  '
$ws.Range("E29").Value = '  -0.20%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0520'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.54%  '
$ws.Range("E31").Value = '  -1.78%  '
$ws.Range("E32").Value = '  +0.18%  '
$ws.Range("E33").Value = '  -2.80%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.79'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.28%  '
$ws.Range("D35").Value = '1.391.75'
$ws.Range("E35").Value = '  -3.54%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.637'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.38%  '
$ws.Range("E37").Value = '  -0.93%  '
$ws.Range("E38").Value = '  -2.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.930'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.42%  '
$ws.Range("E40").Value = '  -0.43%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '78.80'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.24%  '
$ws.Range("E42").Value = '  -4.45%  '
$ws.Range("E43").Value = '  +0.46%  '
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("E45").Value = '  -2.80%  '
$ws.Range("E46").Value = '  -1.56%  '
$ws.Range("D47").Value = '1.926.27'
$ws.Range("E47").Value = '  -0.76%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '104.13'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.02%  '
$ws.Range("E49").Value = '  -0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '11.71'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.37%  '
$ws.Range("E51").Value = '  -2.49%  '
